$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The instruction field (column Y) was binary-encoded numbers; now it should
# explicitly display as hex text. Rows 3-24 become "0x4d", rows 26-57 become "0x7d".
# (Row 25 is an empty spacer row and has no Y value.)

for ($r = 3; $r -le 24; $r++) {
    $ws.Cells.Item($r, 25).Value = "0x4d"
}

for ($r = 26; $r -le 57; $r++) {
    $ws.Cells.Item($r, 25).Value = "0x7d"
}

# Update the saved view state to match the latest interaction position.
$ws.Range("Y60").Select()
